$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.014400005340576
$ws.Range("B1").Value = 2.288256645202637
$ws.Range("C1").Value = 4.823341846466064
$ws.Range("D1").Value = 1.56379759311676
$ws.Range("E1").Value = 1.277548551559448
